# Updates the cryptocurrency tracking data on the active worksheet with the
# latest price/volume snapshot. Two pairs of rows swap their relative
# ranking (rows 12/13: Solana <-> BinanceUSD; rows 50/51: EOS <-> Flow),
# so those rows' Coin/Link/Price/Volume cells are fully rewritten while
# the rest only refresh Price and/or Volume(1h).
#
# Numeric-looking price strings (e.g. "1.005") are entered with a leading
# apostrophe so Excel keeps them as text, matching the column's existing
# plain-text values (prices such as "27.878.95" already aren't parsed as
# numbers, so they don't need the apostrophe).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.878.95'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.751.83'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").Value = "'335.29"
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").Value = "'0.3405"
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = "'45.93"
$ws.Range("E9").Value = '  -2.24%  '
$ws.Range("D10").Value = "'1.114"
$ws.Range("E10").Value = '  -2.44%  '
$ws.Range("D11").Value = "'0.07223"
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = "'22.49"
$ws.Range("E13").Value = '  +1.58%  '
$ws.Range("D14").Value = "'6.163"
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").Value = "'7.142"
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '1.754.04'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("D18").Value = "'0.06609"
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("D19").Value = "'78.88"
$ws.Range("E19").Value = '  -3.88%  '
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = "'16.70"
$ws.Range("E21").Value = '  -3.85%  '
$ws.Range("D22").Value = "'6.224"
$ws.Range("E22").Value = '  -3.36%  '
$ws.Range("D23").Value = '27.904.95'
$ws.Range("E23").Value = '  +1.38%  '
$ws.Range("D24").Value = "'11.66"
$ws.Range("E24").Value = '  -4.17%  '
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = "'153.31"
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").Value = "'19.85"
$ws.Range("E27").Value = '  -4.16%  '
$ws.Range("D28").Value = "'2.311"
$ws.Range("E28").Value = '  -5.21%  '
$ws.Range("D29").Value = '1.954.63'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").Value = "'1.267"
$ws.Range("E30").Value = '  -11.56%  '
$ws.Range("D31").Value = "'132.28"
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("D32").Value = "'4.023"
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("D33").Value = "'5.837"
$ws.Range("E33").Value = '  -4.61%  '
$ws.Range("D34").Value = "'0.08817"
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").Value = "'12.21"
$ws.Range("E35").Value = '  -3.55%  '
$ws.Range("D36").Value = "'0.6571"
$ws.Range("E36").Value = '  -3.05%  '
$ws.Range("D37").Value = "'0.02284"
$ws.Range("E37").Value = '  -6.14%  '
$ws.Range("D38").Value = "'5.143"
$ws.Range("E38").Value = '  -4.45%  '
$ws.Range("D39").Value = "'0.06154"
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").Value = "'1.501"
$ws.Range("E40").Value = '  -2.45%  '
$ws.Range("D41").Value = "'0.2102"
$ws.Range("E41").Value = '  -4.15%  '
$ws.Range("D42").Value = "'1.208"
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("D43").Value = "'7.991"
$ws.Range("E43").Value = '  -5.02%  '
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("D45").Value = "'13.71"
$ws.Range("E45").Value = '  -3.72%  '
$ws.Range("D46").Value = "'0.6069"
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").Value = "'126.51"
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("D49").Value = "'2.003"
$ws.Range("E49").Value = '  -4.77%  '
$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").Value = "'1.121"
$ws.Range("E50").Value = '  +5.74%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = "'1.172"
$ws.Range("E51").Value = '  +2.53%  '
